# Apply updated cryptocurrency price/volume data to the active worksheet
# (values scraped on Tue Jun 25 11:24:00 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.372.25"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.382.52"
$ws.Range("E3").Value = "  +1.97%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.87"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.25"
$ws.Range("E6").Value = "  +9.13%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.381.82"
$ws.Range("E8").Value = "  +1.98%  "
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.61"
$ws.Range("E10").Value = "  +4.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.124"
$ws.Range("E11").Value = "  +5.52%  "
$ws.Range("E12").Value = "  +5.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.956.20"
$ws.Range("E13").Value = "  +1.91%  "
$ws.Range("E14").Value = "  +2.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000174"
$ws.Range("E15").Value = "  +4.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.364.19"
$ws.Range("E16").Value = "  +1.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.23"
$ws.Range("E17").Value = "  +3.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.442.54"
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.96"
$ws.Range("E19").Value = "  +6.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.49"
$ws.Range("E20").Value = "  +4.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.82"
$ws.Range("E21").Value = "  +4.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "382.14"
$ws.Range("E22").Value = "  +9.41%  "
$ws.Range("E23").Value = "  +4.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.515.60"
$ws.Range("E24").Value = "  +1.94%  "
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.89"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000120"
$ws.Range("E27").Value = "  +12.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.68"
$ws.Range("E28").Value = "  +16.44%  "
$ws.Range("E29").Value = "  +9.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.17"
$ws.Range("E31").Value = "  +4.59%  "
$ws.Range("E32").Value = "  +6.35%  "
$ws.Range("E33").Value = "  +1.81%  "
$ws.Range("E35").Value = "  +1.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.49"
$ws.Range("E36").Value = "  +5.22%  "
$ws.Range("E37").Value = "  +3.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.00"
$ws.Range("E38").Value = "  +4.97%  "
$ws.Range("E39").Value = "  +6.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "163.11"
$ws.Range("E40").Value = "  +0.25%  "
$ws.Range("E41").Value = "  +7.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.44"
$ws.Range("E43").Value = "  +5.62%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.54"
$ws.Range("E44").Value = "  +1.54%  "
$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.21"
$ws.Range("E45").Value = "  +9.52%  "
$ws.Range("E46").Value = "  +2.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.68"
$ws.Range("E47").Value = "  +9.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.38"
$ws.Range("E48").Value = "  +4.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.96"
$ws.Range("E49").Value = "  +5.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.23"
$ws.Range("E50").Value = "  +12.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.42"
$ws.Range("E51").Value = "  +13.10%  "
